# Run csv export tool
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSfVBP")

# Replace the formula cells D2:M2 with plain literal zero values
$ws.Range("D2:M2").Value = 0

# Update the active selection on the BSfVBP sheet
$ws.Activate()
$ws.Range("C2:N2").Select()
